$d = $word.ActiveDocument

# --- Bullet 1: research question about technologies/approach for data distribution ---
# Old:  "What technologies are most suitable for developing HeardIT?"
# New:  "What is the most suitable approach for distributing and storing the data used by HeardIT?"
$found1 = $d.Content.Find.Execute(
    "What technologies are most suitable for developing HeardIT?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "What is the most suitable approach for distributing and storing the data used by HeardIT?",
    2)

# --- Bullet 2: research question about deploying HeardIT ---
# Old:  "What technologies and methods are most suitable for deploying HeardIT?"
# New:  "What technologies and methods are most suitable for deploying HeardIT to the cloud?"
$found2 = $d.Content.Find.Execute(
    "What technologies and methods are most suitable for deploying HeardIT?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "What technologies and methods are most suitable for deploying HeardIT to the cloud?",
    2)

Write-Output "Bullet1 replaced: $found1"
Write-Output "Bullet2 replaced: $found2"
